$d = $word.ActiveDocument

# --- Update the date title paragraph (first paragraph, before the table) ---
$d.Content.Find.Execute("2026-02-28 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-03-01 Sunday", 2) | Out-Null

# --- Update each multiplication-fact cell in the table ---
# Cells are addressed by exact (row, col) rather than text search, since a couple of
# old/new values collide across cells (e.g. "18x64=1152" is both a target value in one
# cell and the original value of another), which would make a plain Find/Replace pass
# order-dependent and unsafe.
$tbl = $d.Tables.Item(1)

$cell = $tbl.Cell(1, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "61×39=2379") {
    throw "Unexpected existing text in cell (1,1): $($cell.Range.Text)"
}
$cell.Range.Text = "40×91=3640"

$cell = $tbl.Cell(1, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "21×94=1974") {
    throw "Unexpected existing text in cell (1,2): $($cell.Range.Text)"
}
$cell.Range.Text = "81×12=972"

$cell = $tbl.Cell(1, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "11×48=528") {
    throw "Unexpected existing text in cell (1,3): $($cell.Range.Text)"
}
$cell.Range.Text = "69×42=2898"

$cell = $tbl.Cell(1, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "85×62=5270") {
    throw "Unexpected existing text in cell (1,4): $($cell.Range.Text)"
}
$cell.Range.Text = "85×30=2550"

$cell = $tbl.Cell(1, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "58×93=5394") {
    throw "Unexpected existing text in cell (1,5): $($cell.Range.Text)"
}
$cell.Range.Text = "42×65=2730"

$cell = $tbl.Cell(5, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "75×17=1275") {
    throw "Unexpected existing text in cell (5,1): $($cell.Range.Text)"
}
$cell.Range.Text = "20×42=840"

$cell = $tbl.Cell(5, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "50×27=1350") {
    throw "Unexpected existing text in cell (5,2): $($cell.Range.Text)"
}
$cell.Range.Text = "65×59=3835"

$cell = $tbl.Cell(5, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "65×97=6305") {
    throw "Unexpected existing text in cell (5,3): $($cell.Range.Text)"
}
$cell.Range.Text = "93×85=7905"

$cell = $tbl.Cell(5, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "81×94=7614") {
    throw "Unexpected existing text in cell (5,4): $($cell.Range.Text)"
}
$cell.Range.Text = "86×47=4042"

$cell = $tbl.Cell(5, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "12×72=864") {
    throw "Unexpected existing text in cell (5,5): $($cell.Range.Text)"
}
$cell.Range.Text = "18×64=1152"

$cell = $tbl.Cell(10, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "93×16=1488") {
    throw "Unexpected existing text in cell (10,1): $($cell.Range.Text)"
}
$cell.Range.Text = "63×95=5985"

$cell = $tbl.Cell(10, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "81×42=3402") {
    throw "Unexpected existing text in cell (10,2): $($cell.Range.Text)"
}
$cell.Range.Text = "70×87=6090"

$cell = $tbl.Cell(10, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "98×33=3234") {
    throw "Unexpected existing text in cell (10,3): $($cell.Range.Text)"
}
$cell.Range.Text = "14×51=714"

$cell = $tbl.Cell(10, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "64×18=1152") {
    throw "Unexpected existing text in cell (10,4): $($cell.Range.Text)"
}
$cell.Range.Text = "85×63=5355"

$cell = $tbl.Cell(10, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "27×88=2376") {
    throw "Unexpected existing text in cell (10,5): $($cell.Range.Text)"
}
$cell.Range.Text = "94×84=7896"

$cell = $tbl.Cell(15, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "25×81=2025") {
    throw "Unexpected existing text in cell (15,1): $($cell.Range.Text)"
}
$cell.Range.Text = "53×65=3445"

$cell = $tbl.Cell(15, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "15×28=420") {
    throw "Unexpected existing text in cell (15,2): $($cell.Range.Text)"
}
$cell.Range.Text = "88×80=7040"

$cell = $tbl.Cell(15, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "42×47=1974") {
    throw "Unexpected existing text in cell (15,3): $($cell.Range.Text)"
}
$cell.Range.Text = "50×77=3850"

$cell = $tbl.Cell(15, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "70×89=6230") {
    throw "Unexpected existing text in cell (15,4): $($cell.Range.Text)"
}
$cell.Range.Text = "14×84=1176"

$cell = $tbl.Cell(15, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "64×58=3712") {
    throw "Unexpected existing text in cell (15,5): $($cell.Range.Text)"
}
$cell.Range.Text = "83×44=3652"

$cell = $tbl.Cell(20, 1)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "98×75=7350") {
    throw "Unexpected existing text in cell (20,1): $($cell.Range.Text)"
}
$cell.Range.Text = "36×74=2664"

$cell = $tbl.Cell(20, 2)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "48×45=2160") {
    throw "Unexpected existing text in cell (20,2): $($cell.Range.Text)"
}
$cell.Range.Text = "23×73=1679"

$cell = $tbl.Cell(20, 3)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "91×14=1274") {
    throw "Unexpected existing text in cell (20,3): $($cell.Range.Text)"
}
$cell.Range.Text = "42×53=2226"

$cell = $tbl.Cell(20, 4)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "18×64=1152") {
    throw "Unexpected existing text in cell (20,4): $($cell.Range.Text)"
}
$cell.Range.Text = "48×49=2352"

$cell = $tbl.Cell(20, 5)
if ($cell.Range.Text.TrimEnd([char]7, [char]13) -ne "62×89=5518") {
    throw "Unexpected existing text in cell (20,5): $($cell.Range.Text)"
}
$cell.Range.Text = "69×20=1380"
